$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at Q:R, shifting old Q->S, R->T, S->U, T->V
$ws.Range("Q1:R1").EntireColumn.Insert()

# Header row
$ws.Range("Q1").Value = "default_count"
$ws.Range("R1").Value = "default_value"

# Data rows: set default_count (Q) = 0, default_value (R) = "<Unspecified>", and update the
# most_frequent_value (now in column S) with the new values
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = "<Unspecified>"
$ws.Range("S2").Value = "KELLY"

$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = "<Unspecified>"
$ws.Range("S3").Value = "John"

$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = "<Unspecified>"
$ws.Range("S4").Value = "n/a"

$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = "<Unspecified>"
$ws.Range("S5").Value = "City of Agra"

$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = "<Unspecified>"
$ws.Range("S6").NumberFormat = "@"
$ws.Range("S6").Value = "1878"
$ws.Range("S6").Style = "Normal"

$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = "<Unspecified>"
$ws.Range("S7").NumberFormat = "@"
$ws.Range("S7").Value = "18418"
$ws.Range("S7").Style = "Normal"

$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = "<Unspecified>"
$ws.Range("S8").Value = "IMA3/18"

$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = "<Unspecified>"
$ws.Range("S9").Value = "Z6612"
